$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Row 150: append Discount note to description, mark Done = TRUE
$ws.Range("C150").Value = "#1: استيكر مقاس 5.5*13.5 سم - طباعة ديجيتال | Qty: 100 | Price: 1.03 | Total: 103 | VAT: 0; #2: سليف علبة - طباعة ورق 150 جرام - مقاس 21*6.5 سم | Qty: 100 | Price: 2.07 | Total: 207 | VAT: 0; #3: تسوية الضريبة | Qty: 1 | Price: 0.5 | Total: 0.5 | VAT: 0; Discount: 0"
$ws.Range("I150").Value = $true

# Row 151: append Discount note to description, mark Done = TRUE
$ws.Range("C151").Value = "#1: توريد وطباعة أكياس ورقية - مقاس A5 | Qty: 100 | Price: 3.45 | Total: 345 | VAT: 0; #2: استيكر مقاس 5*5 سم - دائري | Qty: 100 | Price: .69 | Total: 69 | VAT: 0; Discount: 0"
$ws.Range("I151").Value = $true

# Row 155: corrected unit prices, discount note, reference number, updated amount
$ws.Range("C155").Value = "#1: استيكر الأدهم - شفاف - مقاس 3*2 سم | Qty: 100 | Price: 0.747 | Total: 74.7 | VAT: 0; #2: استيكر الكميت - شفاف - مقاس 3*2 سم | Qty: 100 | Price: 0.748 | Total: 74.8 | VAT: 0; #3: استيكر مقاس 5.5*5.5 سم - دائري | Qty: 200 | Price: 0.805 | Total: 161 | VAT: 0; #4: كرت طريقة الاستخدام - كوشيه مسلفن مطفي - مقاس A6 - طباعة وجهين | Qty: 500 | Price: 0.506 | Total: 253 | VAT: 0; Discount: 13.5"
$ws.Range("D155").Value = "QB#5168"
$ws.Range("E155").Value = 563.5
$ws.Range("H155").ClearContents()

# Row 156: corrected price/total/VAT, discount note, reference numbers, mark Done = TRUE
$ws.Range("C156").Value = "#1: ختم دائري مقاس R538 | Qty: 1 | Price: 100 | Total: 100 | VAT: 15; Discount: 0"
$ws.Range("D156").Value = "Invoice#365, QB#5169"
$ws.Range("E156").Value = 100
$ws.Range("F156").Value = 15
$ws.Range("I156").Value = $true

# Row 168: corrected unit prices/totals, discount note, mark Done = TRUE
$ws.Range("C168").Value = "#1: استيكر مقاس 60*12 سم | Qty: 2 | Price: 30 | Total: 60 | VAT: 0; #2: استيكر مقاس 44*7 سم | Qty: 1 | Price: 20 | Total: 20 | VAT: 0; Discount: 0"
$ws.Range("E168").Value = 80
$ws.Range("G168").Value = 80
$ws.Range("I168").Value = $true

# New transaction rows 203-206
$ws.Range("A203").Value = "دوت وان كافيه"
$ws.Range("B203").Value = "'" + "2025-08-21"
$ws.Range("C203").Value = "#1: رول أب على خامة بنر مقاس 85*200 سم | Qty: 3 | Price: 200 | Total: 600 | VAT: 90; Discount: 0"
$ws.Range("E203").Value = 600
$ws.Range("F203").Value = 90
$ws.Range("G203").Value = 690
$ws.Range("I203").Value = $false

$ws.Range("A204").Value = "أبراج اللؤلؤة"
$ws.Range("B204").Value = "'" + "2025-08-23"
$ws.Range("C204").Value = "#1: ختم دائري R538 | Qty: 1 | Price: 120 | Total: 120 | VAT: 18; Discount: 0"
$ws.Range("E204").Value = 120
$ws.Range("F204").Value = 18
$ws.Range("G204").Value = 138
$ws.Range("I204").Value = $false

$ws.Range("A205").Value = "خالد أبو سعيد"
$ws.Range("B205").Value = "'" + "2025-08-23"
$ws.Range("C205").Value = "#1: استيكرات المنتجات | Qty: 1 | Price:  | Total: 0 | VAT: 0; Discount: 0"
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("I205").Value = $false

$ws.Range("A206").Value = "أسامة الأحمدي"
$ws.Range("B206").Value = "'" + "2025-08-23"
$ws.Range("C206").Value = "#1: رول أب مع مكينة - طباعة بنر - مقاس 85*200 سم | Qty: 1 | Price: 170 | Total: 170 | VAT: 0; #2: بنر مقاس 80*140 سم | Qty: 1 | Price: 50 | Total: 50 | VAT: 0; #3: بنر مقاس 120*160 سم | Qty: 1 | Price: 70 | Total: 70 | VAT: 0; Discount: 0"
$ws.Range("E206").Value = 290
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 290
$ws.Range("I206").Value = $false

